$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row total correct count: 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row total marks: 51 -> 85
$ws.Range("B12").Value = 85

# Update the Correct/Total marks text: "50/84" -> "85/140"
$ws.Range("E12").Value = "85/140"
